# "a few more fixes" - apply the four changes captured in the target diff:
#   1. Nudge the x-offset of grouped shape "Group 13" (id 255).
#   2. Nudge the y-offset of "Straight Connector 494" (id 495).
#   3. Remove the stray/duplicate "Straight Connector 375" (id 376).
#   4. Add a new small Oval+Rectangle "Group 13" (duplicated from the
#      existing dot-marker group) positioned near the R/2 resistor.

function Find-ShapeById {
    param($shapes, $targetId)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $cand = $shapes.Item($i)
        if ($cand.Id -eq $targetId) {
            return $cand
        }
    }
    return $null
}

function Find-ShapeIndexById {
    param($shapes, $targetId)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $cand = $shapes.Item($i)
        if ($cand.Id -eq $targetId) {
            return $i
        }
    }
    return 0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Move group "Group 13" (id 255) slightly left -----------------
$grp255 = Find-ShapeById $s.Shapes 255
# Target EMU x=3470313 (was 3483565); nudge chosen so the lossy Single
# round-trip used by Shape.Left lands exactly on the target EMU value.
$grp255.Left = 273.25301212598424

# --- 2. Move connector "Straight Connector 494" (id 495) down --------
$cxn495 = Find-ShapeById $s.Shapes 495
# Target EMU y=1503724 (was 1490472).
$cxn495.Top = 118.40346456692913

# --- 3. Delete the stray connector "Straight Connector 375" (id 376) -
$cxn376 = Find-ShapeById $s.Shapes 376
$cxn376.Delete()

# --- 4. Add a new dot-marker group near the R/2 resistor --------------
# Build it the same way the original author would have: duplicate the
# existing standalone "Oval 14" / "Rectangle 15" marker shapes (which
# already carry the right fill/line formatting) and group them.
$ovalSrc = Find-ShapeById $s.Shapes 387
$rectSrc = Find-ShapeById $s.Shapes 386

$ovalDup = $ovalSrc.Duplicate().Item(1)
$rectDup = $rectSrc.Duplicate().Item(1)

# Absolute target positions (EMU / 12700 = points):
#   oval off = (1327153, 2247106)  ext = (76200, 76200)
#   rect off = (1289053, 2209006)  ext = (152400, 152400)
$ovalDup.Left = 104.50023622047244
$ovalDup.Top = 176.93748031496062
$rectDup.Left = 101.50023622047244
$rectDup.Top = 173.93748031496062

$ovalDupIdx = Find-ShapeIndexById $s.Shapes $ovalDup.Id
$rectDupIdx = Find-ShapeIndexById $s.Shapes $rectDup.Id
$newRange = $s.Shapes.Range(@($ovalDupIdx, $rectDupIdx))
$newGroup = $newRange.Group()
$newGroup.Name = "Group 13"
